# Applies the "handles float input" marksheet re-grade edit:
#  - recompute the summary block (rows 10-12)
#  - drop the third "Student Ans / Correct Ans" column pair (G:H)
#  - drop the now-unused part of the second pair (D:E) except rows 17-18
#  - fill in the student's actual answers in column A (rows 18-40ish),
#    colour-coded green (correct) / red (incorrect) via existing cell styles

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10-12: give column A the same "mtitleStyle" look as the other
#     section headers (re-use the existing style via copy/paste-format so
#     we don't fork a brand new cellXf entry) ---
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 10: Right / Wrong / Not Attempt / Max
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 28

# Row 11: Marking scheme (+4 / -1), now a genuine number instead of text
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Totals
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "54/112"

# --- Drop the third Student/Correct-Ans column pair (G:H) entirely ---
$ws.Range("G15:H21").Clear()

# --- Drop the now-unused tail of the second column pair (D:E), keeping
#     only the header row (15) and the two rows that are still used
#     (17-18) ---
$ws.Range("D19:E40").Clear()

# --- Fill in column D for the two surviving second-pair rows ---
$ws.Range("D15").Copy()
$ws.Range("D17:D18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option B"
# D17 is a correct answer (matches E17) -> correctStyle; D18 is wrong -> incorrectStyle
$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Column A: student answers, colour-coded by correctness ---
# Correct answers (style matches column B -> "correctStyle", green)
$correctRows = @(18, 19, 20, 21, 22, 25, 30, 32, 33, 35, 36, 38, 39)
foreach ($r in $correctRows) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 2).Value2
}
$ws.Range("B10").Copy()
foreach ($r in $correctRows) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false
foreach ($r in $correctRows) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 2).Value2
}

# Incorrect answer (style matches column C -> "incorrectStyle", red)
$ws.Range("C10").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A34").Value = "Option A"

# --- Used range shrinks from A5:H40 to A5:E40 now that G:H is gone ---
$ws.UsedRange | Out-Null
